$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value (written as literal text, matching
# the workbook's existing convention of storing these ticker figures
# as text rather than numbers/percentages).
$updates = @{
    "D2" = "328.11"
    "E2" = "3.60%"
    "D3" = "40.41"
    "E3" = "6.67%"
    "D4" = "5.829"
    "E4" = "12.97%"
    "D5" = "0.08074"
    "E5" = "1.38%"
    "D6" = "4.593"
    "E6" = "2.88%"
    "D7" = "8.770"
    "E7" = "3.39%"
    "D8" = "1.952"
    "E8" = "1.57%"
    "D9" = "2.941"
    "E9" = "-1.29%"
    "D10" = "0.9432"
    "E10" = "0.20%"
    "D11" = "0.1290"
    "E11" = "1.04%"
    "E12" = "1.87%"
    "D13" = "8.933"
    "E13" = "36.31%"
    "D14" = "0.09210"
    "E14" = "2.23%"
    "D15" = "0.03503"
    "E15" = "2.03%"
    "E16" = "0.86%"
    "D17" = "0.001316"
    "E17" = "-5.25%"
    "D18" = "0.006200"
    "E18" = "1.43%"
    "D19" = "3.367"
    "E19" = "-1.57%"
    "D20" = "0.3564"
    "E20" = "1.46%"
    "D21" = "0.1412"
    "E21" = "8.23%"
    "D22" = "0.2413"
    "E22" = "4.89%"
    "E23" = "1.26%"
    "D24" = "0.001261"
    "E24" = "4.75%"
    "D25" = "0.004352"
    "E25" = "-1.32%"
    "D26" = "0.0001144"
    "E26" = "-13.63%"
    "E27" = "0.45%"
    "D39" = "0.02418"
    "E39" = "1.31%"
    "D40" = "0.05278"
    "E40" = "2.04%"
    "D41" = "0.007479"
    "E41" = "0.72%"
    "D42" = "0.1428"
    "E42" = "2.29%"
    "D43" = "0.008697"
    "E43" = "3.16%"
    "D44" = "0.002109"
    "E44" = "0.33%"
    "D45" = "0.01098"
    "E45" = "25.54%"
    "D46" = "0.00006895"
    "E46" = "6.55%"
    "D47" = "0.00000000753"
    "E47" = "0.76%"
    "D48" = "0.003166"
    "E48" = "10.84%"
    "D49" = "0.001703"
    "E49" = "1.18%"
    "D50" = "0.00002108"
    "E50" = "0.76%"
    "D51" = "0.0002008"
    "E51" = "0.76%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    # Leading apostrophe forces Excel to store the entry as literal text
    # (these columns hold price/volume figures as strings, not numbers).
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = $origStyle
}
